$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.264243666666667
$ws.Range("H2").Value = 6.792731
$ws.Range("I2").Value = 0.4114976873616865
$ws.Range("J2").Value = 0.4114976873616865
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.922246666666667
$ws.Range("N2").Value = 5.76674
$ws.Range("O2").Value = 0.04282684507083529
$ws.Range("P2").Value = 0.04282684507083529
$ws.Range("Q2").Value = 4.352434840771112
$ws.Range("R2").Value = 39.17191356694
$ws.Range("S2").Value = 0.01762314770364597
$ws.Range("T2").Value = 0.01762314770364597
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.264243666666667
$ws.Range("H3").Value = 6.792731
$ws.Range("I3").Value = 0.4114976873616865
$ws.Range("J3").Value = 0.4114976873616865
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 21.84955866666667
$ws.Range("N3").Value = 65.548676
$ws.Range("O3").Value = 0.4867989525538483
$ws.Range("P3").Value = 0.4867989525538483
$ws.Range("Q3").Value = 49.47272483046178
$ws.Range("R3").Value = 445.254523474156
$ws.Range("S3").Value = 0.2003166431859999
$ws.Range("T3").Value = 0.2003166431859999
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.264243666666667
$ws.Range("H4").Value = 6.792731
$ws.Range("I4").Value = 0.4114976873616865
$ws.Range("J4").Value = 0.4114976873616865
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.102188
$ws.Range("N4").Value = 0.306564
$ws.Range("O4").Value = 0.002276705544605019
$ws.Range("P4").Value = 0.002276705544605019
$ws.Range("Q4").Value = 0.2313785318093333
$ws.Range("R4").Value = 2.082406786284
$ws.Range("S4").Value = 0.0009368590664084943
$ws.Range("T4").Value = 0.0009368590664084943
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.264243666666667
$ws.Range("H5").Value = 6.792731
$ws.Range("I5").Value = 0.4114976873616865
$ws.Range("J5").Value = 0.4114976873616865
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.01015966666667
$ws.Range("N5").Value = 63.030479
$ws.Range("O5").Value = 0.4680974968307114
$ws.Range("P5").Value = 0.4680974968307114
$ws.Range("Q5").Value = 47.57212096090544
$ws.Range("R5").Value = 428.149088648149
$ws.Range("S5").Value = 0.1926210374056321
$ws.Range("T5").Value = 0.1926210374056321
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.089228666666667
$ws.Range("H6").Value = 6.267686
$ws.Range("I6").Value = 0.3796909216792509
$ws.Range("J6").Value = 0.3796909216792509
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.922246666666667
$ws.Range("N6").Value = 5.76674
$ws.Range("O6").Value = 0.04282684507083529
$ws.Range("P6").Value = 0.04282684507083529
$ws.Range("Q6").Value = 4.016012840404445
$ws.Range("R6").Value = 36.14411556364001
$ws.Range("S6").Value = 0.01626096427755994
$ws.Range("T6").Value = 0.01626096427755994
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.089228666666667
$ws.Range("H7").Value = 6.267686
$ws.Range("I7").Value = 0.3796909216792509
$ws.Range("J7").Value = 0.3796909216792509
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.84955866666667
$ws.Range("N7").Value = 65.548676
$ws.Range("O7").Value = 0.4867989525538483
$ws.Range("P7").Value = 0.4867989525538483
$ws.Range("Q7").Value = 45.64872432041511
$ws.Range("R7").Value = 410.838518883736
$ws.Range("S7").Value = 0.1848331429676646
$ws.Range("T7").Value = 0.1848331429676646
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.089228666666667
$ws.Range("H8").Value = 6.267686
$ws.Range("I8").Value = 0.3796909216792509
$ws.Range("J8").Value = 0.3796909216792509
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.102188
$ws.Range("N8").Value = 0.306564
$ws.Range("O8").Value = 0.002276705544605019
$ws.Range("P8").Value = 0.002276705544605019
$ws.Range("Q8").Value = 0.2134940989893333
$ws.Range("R8").Value = 1.921446890904
$ws.Range("S8").Value = 0.0008644444266233405
$ws.Range("T8").Value = 0.0008644444266233405
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.089228666666667
$ws.Range("H9").Value = 6.267686
$ws.Range("I9").Value = 0.3796909216792509
$ws.Range("J9").Value = 0.3796909216792509
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 21.01015966666667
$ws.Range("N9").Value = 63.030479
$ws.Range("O9").Value = 0.4680974968307114
$ws.Range("P9").Value = 0.4680974968307114
$ws.Range("Q9").Value = 43.89502786684378
$ws.Range("R9").Value = 395.055250801594
$ws.Range("S9").Value = 0.1777323700074031
$ws.Range("T9").Value = 0.1777323700074031
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1905406666666667
$ws.Range("H10").Value = 0.571622
$ws.Range("I10").Value = 0.03462835949856721
$ws.Range("J10").Value = 0.03462835949856721
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.922246666666667
$ws.Range("N10").Value = 5.76674
$ws.Range("O10").Value = 0.04282684507083529
$ws.Range("P10").Value = 0.04282684507083529
$ws.Range("Q10").Value = 0.3662661613644445
$ws.Range("R10").Value = 3.29639545228
$ws.Range("S10").Value = 0.001483023387302326
$ws.Range("T10").Value = 0.001483023387302326
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1905406666666667
$ws.Range("H11").Value = 0.571622
$ws.Range("I11").Value = 0.03462835949856721
$ws.Range("J11").Value = 0.03462835949856721
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 21.84955866666667
$ws.Range("N11").Value = 65.548676
$ws.Range("O11").Value = 0.4867989525538483
$ws.Range("P11").Value = 0.4867989525538483
$ws.Range("Q11").Value = 4.163229474719111
$ws.Range("R11").Value = 37.469065272472
$ws.Range("S11").Value = 0.01685704913256062
$ws.Range("T11").Value = 0.01685704913256062
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1905406666666667
$ws.Range("H12").Value = 0.571622
$ws.Range("I12").Value = 0.03462835949856721
$ws.Range("J12").Value = 0.03462835949856721
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.102188
$ws.Range("N12").Value = 0.306564
$ws.Range("O12").Value = 0.002276705544605019
$ws.Range("P12").Value = 0.002276705544605019
$ws.Range("Q12").Value = 0.01947096964533333
$ws.Range("R12").Value = 0.175238726808
$ws.Range("S12").Value = 0.00007883857807096385
$ws.Range("T12").Value = 0.00007883857807096385
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1905406666666667
$ws.Range("H13").Value = 0.571622
$ws.Range("I13").Value = 0.03462835949856721
$ws.Range("J13").Value = 0.03462835949856721
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 21.01015966666667
$ws.Range("N13").Value = 63.030479
$ws.Range("O13").Value = 0.4680974968307114
$ws.Range("P13").Value = 0.4680974968307114
$ws.Range("Q13").Value = 4.003289829659778
$ws.Range("R13").Value = 36.029608466938
$ws.Range("S13").Value = 0.0162094484006333
$ws.Range("T13").Value = 0.0162094484006333
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.9584326666666668
$ws.Range("H14").Value = 2.875298
$ws.Range("I14").Value = 0.1741830314604954
$ws.Range("J14").Value = 0.1741830314604954
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.922246666666667
$ws.Range("N14").Value = 5.76674
$ws.Range("O14").Value = 0.04282684507083529
$ws.Range("P14").Value = 0.04282684507083529
$ws.Range("Q14").Value = 1.842343998724445
$ws.Range("R14").Value = 16.58109598852
$ws.Range("S14").Value = 0.007459709702327067
$ws.Range("T14").Value = 0.007459709702327067
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.9584326666666668
$ws.Range("H15").Value = 2.875298
$ws.Range("I15").Value = 0.1741830314604954
$ws.Range("J15").Value = 0.1741830314604954
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 21.84955866666667
$ws.Range("N15").Value = 65.548676
$ws.Range("O15").Value = 0.4867989525538483
$ws.Range("P15").Value = 0.4867989525538483
$ws.Range("Q15").Value = 20.94133077838312
$ws.Range("R15").Value = 188.471977005448
$ws.Range("S15").Value = 0.08479211726762319
$ws.Range("T15").Value = 0.08479211726762319
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.9584326666666668
$ws.Range("H16").Value = 2.875298
$ws.Range("I16").Value = 0.1741830314604954
$ws.Range("J16").Value = 0.1741830314604954
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.102188
$ws.Range("N16").Value = 0.306564
$ws.Range("O16").Value = 0.002276705544605019
$ws.Range("P16").Value = 0.002276705544605019
$ws.Range("Q16").Value = 0.09794031734133335
$ws.Range("R16").Value = 0.8814628560720001
$ws.Range("S16").Value = 0.0003965634735022204
$ws.Range("T16").Value = 0.0003965634735022204
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.9584326666666668
$ws.Range("H17").Value = 2.875298
$ws.Range("I17").Value = 0.1741830314604954
$ws.Range("J17").Value = 0.1741830314604954
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 21.01015966666667
$ws.Range("N17").Value = 63.030479
$ws.Range("O17").Value = 0.4680974968307114
$ws.Range("P17").Value = 0.4680974968307114
$ws.Range("Q17").Value = 20.13682335641578
$ws.Range("R17").Value = 181.231410207742
$ws.Range("S17").Value = 0.08153464101704297
$ws.Range("T17").Value = 0.08153464101704297
$wb.Save()
